# Saldo.xlsx update:
#  - add MARCELO (000772433 / 5000) just above MONICA (004387250)
#  - remove CRISTINA (004368994 / 1773.32)
#  - add WANDER (004216504 / 197.18) just above LOHRAN (004332747)
#  - add ILTON (004211368 / 158.5) just above RAFAEL (004212409)
#  - remove the old WANDER row (004216504 / 0.43) further down the sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$col = $ws.Columns.Item(1)

# 1) Insert MARCELO right before the row that holds account 004387250
$anchor = $col.Find("004387250")
$r = $anchor.Row
$ws.Rows.Item($r).EntireRow.Insert()
$ws.Range("A" + $r).NumberFormat = "@"
$ws.Range("A" + $r).Value = "000772433"
$ws.Range("B" + $r).Value = "MARCELO"
$ws.Range("C" + $r).Value = 5000

# 2) Delete the CRISTINA row (004368994 / 1773.32)
$cristina = $col.Find("004368994")
$ws.Rows.Item($cristina.Row).EntireRow.Delete()

# 3) Insert WANDER (004216504 / 197.18) right before LOHRAN (004332747)
$anchor = $col.Find("004332747")
$r = $anchor.Row
$ws.Rows.Item($r).EntireRow.Insert()
$ws.Range("A" + $r).NumberFormat = "@"
$ws.Range("A" + $r).Value = "004216504"
$ws.Range("B" + $r).Value = "WANDER"
$ws.Range("C" + $r).Value = 197.18

# 4) Insert ILTON (004211368 / 158.5) right before RAFAEL (004212409)
$anchor = $col.Find("004212409")
$r = $anchor.Row
$ws.Rows.Item($r).EntireRow.Insert()
$ws.Range("A" + $r).NumberFormat = "@"
$ws.Range("A" + $r).Value = "004211368"
$ws.Range("B" + $r).Value = "ILTON"
$ws.Range("C" + $r).Value = 158.5

# 5) Delete the stale WANDER row further down (004216504 / 0.43).
#    FindNext from the anchor we just inserted to make sure we skip the
#    WANDER row added in step 3 and land on the original lower one.
$first = $col.Find("004216504")
$dup = $col.FindNext($first)
if ($dup.Row -ne $first.Row) {
  $ws.Rows.Item($dup.Row).EntireRow.Delete()
} else {
  $ws.Rows.Item($first.Row).EntireRow.Delete()
}
